# Add a "Neuron_type" column (all "PYR") to Sheet1, inserted before the
# existing "Analyze" (F) / "Session_path" (G) columns, shifting those two
# columns one slot to the right, and shifting the trailing "Comments"
# column (P -> Q) as well.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$xlPasteFormats = -4122
$lastRow = 15

# 1) Shift the "Comments" column (old P1:P9, the only populated rows) to
#    its new location (Q1:Q9), carrying both value and formatting, then
#    clear the old column.
$ws.Range("P1:P9").Copy($ws.Range("Q1:Q9"))
$ws.Range("P1:P9").Clear()

# 2) Shift existing columns G (Session_path) -> H, then F (Analyze) -> G.
#    Do the rightmost column first so the source data for the next step
#    isn't clobbered.
$ws.Range("G1:G15").Copy($ws.Range("H1:H15"))
$ws.Range("F1:F15").Copy($ws.Range("G1:G15"))

# 3) Populate the new column F with the "Neuron_type" header and "PYR"
#    values for every data row, re-using the formatting that now lives in
#    column G (copied from the original F column).
$ws.Range("G1").Copy()
$ws.Range("F1").PasteSpecial($xlPasteFormats)
$ws.Cells.Item(1, 6).Value2 = "Neuron_type"

$ws.Range("G2").Copy()
$ws.Range("F2:F15").PasteSpecial($xlPasteFormats)
for ($r = 2; $r -le $lastRow; $r++) {
    $ws.Cells.Item($r, 6).Value2 = "PYR"
}

$excel.CutCopyMode = 0

# 4) Column widths: the old "Analyze" column width (7.85546875 chars) now
#    belongs to column G, and the new "Neuron_type" column (F) gets its own,
#    wider, width (14.28515625 chars). ColumnWidth uses a slightly different
#    unit than the stored XML "width" (XML width = ColumnWidth + 5/6), so the
#    values below are chosen to land as close as possible on the target.
$ws.Columns.Item(7).ColumnWidth = 7.022135416666667
$ws.Columns.Item(6).ColumnWidth = 13.451822916666666

# 5) Selection, matching the saved view state.
$ws.Range("F2:F15").Select()
